$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.114.87"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "3.150.46"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "3.142.65"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("E10").Value = "  +0.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.56%  "

# Row 12
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("E13").Value = "  -1.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
$ws.Range("D15").Value = "3.672.75"
$ws.Range("E15").Value = "  +0.77%  "

# Row 16
$ws.Range("E16").Value = "  -1.31%  "

# Row 17
$ws.Range("E17").Value = "  +1.52%  "

# Row 18
$ws.Range("D18").Value = "63.943.03"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "3.146.43"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "

# Row 25
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.01%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "

# Row 27
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.18%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.76%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "

# Row 31
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("E32").Value = "  +0.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.94%  "

# Row 34
$ws.Range("E34").Value = "  +2.14%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0839"
$ws.Range("E35").Value = "  -4.50%  "

# Row 36
$ws.Range("E36").Value = "  +1.99%  "

# Row 37
$ws.Range("E37").Value = "  +0.71%  "

# Row 38
$ws.Range("E38").Value = "  -2.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.08%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "

# Row 42
$ws.Range("E42").Value = "  +5.28%  "

# Row 43
$ws.Range("E43").Value = "  +5.07%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0373"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.929.57"
$ws.Range("E45").Value = "  +0.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.108"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.72%  "

# Row 50
$ws.Range("E50").Value = "  +2.57%  "

# Row 51
$ws.Range("E51").Value = "  -0.74%  "
